$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# (matching the original inlineStr/text cell type) and then have their
# formatting cleared so no extra "Text" number format style is introduced.
$numericLookingCells = @(
    'D4',
    'D5',
    'D6',
    'D8',
    'D10',
    'D11',
    'D13',
    'D14',
    'D15',
    'D16',
    'D17',
    'D19',
    'D20',
    'D22',
    'D23',
    'D24',
    'D25',
    'D27',
    'D28',
    'D29',
    'D30',
    'D34',
    'D35',
    'D36',
    'D37',
    'D39',
    'D40',
    'D41',
    'D44',
    'D45',
    'D47',
    'D48',
    'D49',
)

foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values (text form preserved for every cell)
$ws.Range('D2').Value = '29.303.87'
$ws.Range('E2').Value = '  -0.70%  '
$ws.Range('D3').Value = '1.840.08'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '239.04'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('D6').Value = '0.6260'
$ws.Range('E6').Value = '  -0.65%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.07369'
$ws.Range('E8').Value = '  -1.11%  '
$ws.Range('E9').Value = '  -0.58%  '
$ws.Range('D10').Value = '24.73'
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('D11').Value = '0.07720'
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').Value = '1.832.71'
$ws.Range('E12').Value = '  -0.89%  '
$ws.Range('D13').Value = '4.943'
$ws.Range('E13').Value = '  -1.38%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').Value = '0.00001049'
$ws.Range('E14').Value = '  +2.54%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '0.6626'
$ws.Range('E15').Value = '  -2.90%  '
$ws.Range('D16').Value = '81.36'
$ws.Range('E16').Value = '  -1.45%  '
$ws.Range('D17').Value = '6.233'
$ws.Range('E17').Value = '  -1.37%  '
$ws.Range('D18').Value = '29.330.87'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').Value = '235.51'
$ws.Range('E19').Value = '  +2.51%  '
$ws.Range('D20').Value = '12.22'
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').Value = '7.232'
$ws.Range('E22').Value = '  -3.78%  '
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').Value = '157.34'
$ws.Range('E24').Value = '  -1.20%  '
$ws.Range('D25').Value = '8.412'
$ws.Range('E25').Value = '  -1.13%  '
$ws.Range('E26').Value = '  -1.95%  '
$ws.Range('D27').Value = '17.26'
$ws.Range('D28').Value = '0.07079'
$ws.Range('E28').Value = '  +7.58%  '
$ws.Range('D29').Value = '1.468'
$ws.Range('E29').Value = '  +0.30%  '
$ws.Range('D30').Value = '1.479'
$ws.Range('E30').Value = '  -0.69%  '
$ws.Range('E31').Value = '  -1.95%  '
$ws.Range('E32').Value = '  -1.65%  '
$ws.Range('E33').Value = '  +0.91%  '
$ws.Range('D34').Value = '1.781'
$ws.Range('E34').Value = '  -3.61%  '
$ws.Range('D35').Value = '0.6858'
$ws.Range('E35').Value = '  -1.58%  '
$ws.Range('D36').Value = '2.578'
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('D37').Value = '0.01823'
$ws.Range('E37').Value = '  -2.57%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.233.27'
$ws.Range('E38').Value = '  -1.61%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.778'
$ws.Range('E39').Value = '  -1.87%  '
$ws.Range('D40').Value = '6.718'
$ws.Range('E40').Value = '  -0.97%  '
$ws.Range('D41').Value = '0.9435'
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').Value = '1.998.88'
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('D44').Value = '101.17'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').Value = '65.08'
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('E46').Value = '  +6.90%  '
$ws.Range('D47').Value = '6.931'
$ws.Range('E47').Value = '  -2.28%  '
$ws.Range('D48').Value = '1.683'
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('D49').Value = '8.853'
$ws.Range('E49').Value = '  -1.90%  '
$ws.Range('E51').Value = '  -1.44%  '

# Clear number-format overrides on the cells we text-forced above so the
# resulting cells keep the default (no explicit) style, same as before.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).ClearFormats()
}
